$d = $word.ActiveDocument

# Namespace declaration needed for the injected run XML fragments.
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. "Assignment Name - RCNN & YOLO Assignment" paragraph -----------
# Original runs: "Assignment Name " | proofErr(gramStart) | "-  " | "RCNN"
#                | proofErr(gramEnd) | " & YOLO Assignment"
# Target: two runs, "Assignment Name -  " and "RCNN & YOLO Assignment",
# with both proofErr markers removed.
$assignPara = $d.Paragraphs(4)
$assignRange = $d.Range($assignPara.Range.Start, $assignPara.Range.End)
$assignXml = "<w:p $wns>" +
    '<w:r><w:t xml:space="preserve">Assignment Name -  </w:t></w:r>' +
    '<w:r><w:t>RCNN &amp; YOLO Assignment</w:t></w:r>' +
    '</w:p>'
$assignRange.InsertXML($assignXml)

# --- 2. "Git Link -<url>" paragraph: update the notebook URL -----------
# Keep the "Git Link ", "-" and URL runs distinct, only change the URL text.
$gitPara = $d.Paragraphs(6)
$gitRange = $d.Range($gitPara.Range.Start, $gitPara.Range.End)
$gitXml = "<w:p $wns>" +
    '<w:r><w:t xml:space="preserve">Git Link </w:t></w:r>' +
    '<w:r><w:t>-</w:t></w:r>' +
    '<w:r><w:t>https://github.com/ankitsharma5911/deeplearning-assignment/blob/main/RCNN%26Yolo.ipynb</w:t></w:r>' +
    '</w:p>'
$gitRange.InsertXML($gitXml)
